$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''39.048.41'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -4.60%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''2.238.18'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -7.29%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.24%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''295.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -6.00%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''80.37'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -8.87%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.506'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -5.43%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.21%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.457'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -7.42%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.0770'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -7.33%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('B11').Value = '''OKB'
$ws.Range('B11').Style = 'Normal'
$ws.Range('C11').Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('C11').Style = 'Normal'
$ws.Range('D11').Value = '''47.38'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -11.11%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').Value = '''Avalanche'
$ws.Range('B12').Style = 'Normal'
$ws.Range('C12').Value = '''https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').Value = '''27.85'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -10.95%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = '''  -1.28%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''2.585.02'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -7.28%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''6.06'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -10.09%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''14.03'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -9.33%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''2.246.92'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -6.90%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''0.712'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -7.20%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''38.930.91'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -4.71%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''0.0₃0855'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -6.79%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''5.77'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -7.08%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''65.38'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -7.22%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''9.92'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -8.23%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''226.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -4.83%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -0.09%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''2.38'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -10.42%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''1.73'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -5.64%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('B28').Value = '''EthereumClassic'
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').Value = '''22.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -6.92%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('B29').Value = '''Toncoin'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').Value = '''2.18'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -1.97%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''8.81'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  -6.90%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''147.62'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -6.09%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''31.47'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -7.45%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -0.48%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = '''Filecoin'
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = '''4.75'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -9.44%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('B35').Value = '''WEMIXToken'
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = '''2.31'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -6.27%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''0.0681'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -7.21%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = '''  -3.93%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('B38').Value = '''LidoDAOToken'
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = '''2.61'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -8.27%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = '''Kaspa'
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = '''0.0942'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -5.40%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''14.67'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -9.04%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''1.59'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -8.90%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''3.61'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  -5.98%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = '''Maker'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = '''1.915.70'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -3.39%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = '''ApeXProtocol'
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = '''https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = '''2.20'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -3.66%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''0.0252'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -7.29%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''16.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -9.88%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = '''  -4.05%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''2.51'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -11.52%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''2.478.93'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -6.44%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''87.73'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -6.21%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''66.26'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -10.03%  '
$ws.Range('E51').Style = 'Normal'
